# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.795.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.366.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.36%  "

$ws.Range("E7").Value = "  -5.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.358.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.62%  "

$ws.Range("E10").Value = "  -12.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.16%  "

$ws.Range("E13").Value = "  -11.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.903.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "612.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.666.65"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  -3.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.369.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.913"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.89%  "

$ws.Range("E26").Value = "  -10.36%  "

$ws.Range("E27").Value = "  -9.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.05%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.31%  "

$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.764.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "532.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("E39").Value = "  +38.15%  "

$ws.Range("E40").Value = "  -5.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0725"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.352"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.55%  "

$ws.Range("E44").Value = "  -7.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0418"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.71%  "

$ws.Range("E48").Value = "  -12.88%  "

$ws.Range("E49").Value = "  -7.82%  "

$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("E51").Value = "  -9.95%  "

